# Insert a new weekly price record for "Choclo" (Dulce o Americano, Arica y
# Parinacota) ahead of the existing row 195, shifting the remaining rows of
# the table down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(195).Insert()

$ws.Cells.Item(195, 1).Value = 9
$ws.Cells.Item(195, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(195, 3).Value = "Metropolitana"
$ws.Cells.Item(195, 4).Value = 44466
$ws.Cells.Item(195, 5).Value = 13
$ws.Cells.Item(195, 6).Value = 100112024
$ws.Cells.Item(195, 7).Value = "Choclo"
$ws.Cells.Item(195, 8).Value = "Dulce o Americano"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 43
$ws.Cells.Item(195, 11).Value = 38000
$ws.Cells.Item(195, 12).Value = 40000
$ws.Cells.Item(195, 13).Value = 39023
$ws.Cells.Item(195, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(195, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(195, 16).Value = 557
$ws.Cells.Item(195, 17).Value = 70
$ws.Cells.Item(195, 18).Value = "Hortaliza"

# Re-apply the date-column number format (yyyy-mm-dd hh:mm:ss) to the new
# row's date cell so it matches the rest of column D.
$ws.Cells.Item(195, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat
